# Auto-generated update of Leve profit-calculation worksheets
# Source data: market-board average prices / computed Leve profits
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 66270
$ws.Range("J70").Value = 177304.75
$ws.Range("L70").Value = 531914.25
$ws.Range("N70").Value = -532454.25

# Row 73
$ws.Range("H73").Value = 66270
$ws.Range("J73").Value = 177304.75
$ws.Range("L73").Value = 531914.25
$ws.Range("N73").Value = -533786.25

# Row 98
$ws.Range("H98").Value = 1697.762
$ws.Range("I98").Value = 1604.2106
$ws.Range("J98").Value = 2586.5
$ws.Range("K98").Value = 1604.2106
$ws.Range("L98").Value = 2586.5
$ws.Range("M98").Value = -106.2106000000001
$ws.Range("N98").Value = -5582.5

# Row 122
$ws.Range("H122").Value = 1697.762
$ws.Range("I122").Value = 1604.2106
$ws.Range("J122").Value = 2586.5
$ws.Range("K122").Value = 4812.6318
$ws.Range("L122").Value = 7759.5
$ws.Range("M122").Value = -2362.6318
$ws.Range("N122").Value = -12659.5

# Row 138
$ws.Range("H138").Value = 7992.0527
$ws.Range("I138").Value = 8752.388999999999
$ws.Range("J138").Value = 7756.0864
$ws.Range("K138").Value = 26257.167
$ws.Range("L138").Value = 23268.2592
$ws.Range("M138").Value = -21117.167
$ws.Range("N138").Value = -33548.2592

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 157.75
$ws.Range("J5").Value = 199.5
$ws.Range("L5").Value = 199.5
$ws.Range("N5").Value = -423.5

# Row 31
$ws.Range("H31").Value = 7000
$ws.Range("I31").Value = 7000
$ws.Range("K31").Value = 7000
$ws.Range("M31").Value = -6706

# Row 32
$ws.Range("H32").Value = 15921.204
$ws.Range("I32").Value = 8536.826999999999
$ws.Range("J32").Value = 30197.666
$ws.Range("K32").Value = 8536.826999999999
$ws.Range("L32").Value = 30197.666
$ws.Range("M32").Value = -8249.826999999999
$ws.Range("N32").Value = -30771.666

# Row 45
$ws.Range("H45").Value = 2375
$ws.Range("I45").Value = 1776.7858
$ws.Range("K45").Value = 1776.7858
$ws.Range("M45").Value = -1399.7858

# Row 97
$ws.Range("H97").Value = 965.58826
$ws.Range("I97").Value = 992.3333
$ws.Range("J97").Value = 765
$ws.Range("K97").Value = 992.3333
$ws.Range("L97").Value = 765
$ws.Range("M97").Value = -496.3333
$ws.Range("N97").Value = -1757

# Row 110
$ws.Range("H110").Value = 3238.3845
$ws.Range("I110").Value = 3238.3845
$ws.Range("K110").Value = 3238.3845
$ws.Range("M110").Value = -1193.3845

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 157.75
$ws.Range("J4").Value = 199.5
$ws.Range("L4").Value = 199.5
$ws.Range("N4").Value = -429.5

# Row 99
$ws.Range("H99").Value = 890.1739
$ws.Range("I99").Value = 867.0526
$ws.Range("K99").Value = 867.0526
$ws.Range("M99").Value = 630.9474

# Row 102
$ws.Range("H102").Value = 9000
$ws.Range("I102").Value = 9000
$ws.Range("K102").Value = 9000
$ws.Range("M102").Value = -5755

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 766.3333
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50

# Row 31
$ws.Range("H31").Value = 5024
$ws.Range("I31").Value = 3652.25
$ws.Range("J31").Value = 6290.231
$ws.Range("K31").Value = 3652.25
$ws.Range("L31").Value = 6290.231
$ws.Range("M31").Value = -3357.25
$ws.Range("N31").Value = -6880.231

# Row 34
$ws.Range("H34").Value = 5024
$ws.Range("I34").Value = 3652.25
$ws.Range("J34").Value = 6290.231
$ws.Range("K34").Value = 3652.25
$ws.Range("L34").Value = 6290.231
$ws.Range("M34").Value = -3450.25
$ws.Range("N34").Value = -6694.231

# Row 68
$ws.Range("H68").Value = 37499.5
$ws.Range("J68").Value = 37499.5
$ws.Range("L68").Value = 37499.5
$ws.Range("N68").Value = -38997.5

# Row 71
$ws.Range("H71").Value = 37499.5
$ws.Range("J71").Value = 37499.5
$ws.Range("L71").Value = 112498.5
$ws.Range("N71").Value = -119986.5

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3006.6667
$ws.Range("I3").Value = 2010
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 6030
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -5918
$ws.Range("N3").Value = -15224

# Row 10
$ws.Range("H10").Value = 28.25
$ws.Range("I10").Value = 28.25
$ws.Range("K10").Value = 84.75
$ws.Range("M10").Value = 54.25

# Row 16
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 4
$ws.Range("K16").Value = 12
$ws.Range("M16").Value = 161

# Row 140
$ws.Range("H140").Value = 5272.75
$ws.Range("I140").Value = 4697
$ws.Range("K140").Value = 14091
$ws.Range("M140").Value = -8911

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1001.13336
$ws.Range("I93").Value = 648.3333
$ws.Range("K93").Value = 648.3333
$ws.Range("M93").Value = 599.6667

# Row 136
$ws.Range("H136").Value = 3598.8333
$ws.Range("I136").Value = 3598.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10796.4999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8246.499899999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4249.615
$ws.Range("J81").Value = 5116.3335
$ws.Range("L81").Value = 10232.667
$ws.Range("N81").Value = -12354.667

# Row 84
$ws.Range("H84").Value = 4249.615
$ws.Range("J84").Value = 5116.3335
$ws.Range("L84").Value = 51163.335
$ws.Range("N84").Value = -61771.335

# Row 107
$ws.Range("H107").Value = 1742.375
$ws.Range("I107").Value = 747.25
$ws.Range("K107").Value = 2241.75
$ws.Range("M107").Value = -321.75

# Row 113
$ws.Range("H113").Value = 1238.4615
$ws.Range("I113").Value = 878.75
$ws.Range("J113").Value = 1546.7858
$ws.Range("K113").Value = 2636.25
$ws.Range("L113").Value = 4640.357400000001
$ws.Range("M113").Value = -466.25
$ws.Range("N113").Value = -8980.357400000001

# Row 132
$ws.Range("H132").Value = 2052.3333
$ws.Range("I132").Value = 548.6
$ws.Range("J132").Value = 3126.4285
$ws.Range("K132").Value = 1645.8
$ws.Range("L132").Value = 9379.2855
$ws.Range("M132").Value = 884.1999999999998
$ws.Range("N132").Value = -14439.2855

# Row 136
$ws.Range("H136").Value = 73883.5
$ws.Range("I136").Value = 1497.7778
$ws.Range("J136").Value = 204177.8
$ws.Range("K136").Value = 4493.3334
$ws.Range("L136").Value = 612533.3999999999
$ws.Range("M136").Value = -1943.3334
$ws.Range("N136").Value = -617633.3999999999
